$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 32: 2019-11-13, 16:00 -> 00:00, interruption 120, delta 360 ---
$ws.Range("A32").Value = 43782
$ws.Range("B32").Value = 0.66666666666666663
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 120
$ws.Range("E32").Value = 360

$f32 = "팀원" + " 결과물 취합 및 SAD 1.0 제작"
$ws.Range("F32").Value = $f32
$ws.Range("F32").Characters(3, 20).Font.Name = "맑은 고딕"
$ws.Range("F32").Characters(3, 20).Font.Size = 10

# --- Row 33: 2019-11-14, 01:00 -> 03:00, interruption 0, delta 120 ---
$ws.Range("A33").Value = 43783
$ws.Range("B33").Value = 0.041666666666666664
$ws.Range("C33").Value = 0.125
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 120

$f33 = "SAD 1.0 발표" + " 준비"
$ws.Range("F33").Value = $f33
$ws.Range("F33").Characters(11, 3).Font.Name = "돋움"
$ws.Range("F33").Characters(11, 3).Font.Size = 10

# --- Update the view: scroll position + active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F34").Select()
